# DEPA_Student_Quickstart.docx edit script
# 1) Title: "DEPA Workshop - Student Quickstart Guide (Updated)"
#    -> "DEPA Workshop - Student Quickstart Guide", split into three runs
#       with a spell-check proofErr wrap around "Quickstart".
# 2) "Detailed instructions ..." paragraph: merge the five separate runs
#    into a single run with the same concatenated text.

$d = $word.ActiveDocument

function Get-PkgXml($innerParagraphXml) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$innerParagraphXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

# ---------------------------------------------------------------------
# Edit 1: title paragraph
# ---------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
# Exclude the trailing paragraph mark so only the run content (not the
# paragraph's own <w:pPr>) gets replaced.
$titleContent = $d.Range($titleRange.Start, $titleRange.End - 1)

if ($titleContent.Text -ne "DEPA Workshop – Student Quickstart Guide (Updated)") {
    throw "Unexpected title text: $($titleContent.Text)"
}

$titleXml = @'
<w:p>
<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">DEPA Workshop &#8211; Student </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>Quickstart</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> Guide</w:t></w:r>
</w:p>
'@

$titleContent.InsertXML((Get-PkgXml $titleXml))

# ---------------------------------------------------------------------
# Edit 2: "Detailed instructions..." paragraph - collapse 5 runs into 1
# ---------------------------------------------------------------------

$detailPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.StartsWith("Detailed instructions on ESP setup")) {
        $detailPara = $candidate
        break
    }
}
if ($null -eq $detailPara) {
    throw "Could not locate 'Detailed instructions' paragraph"
}

$detailRange = $detailPara.Range
$detailContent = $d.Range($detailRange.Start, $detailRange.End - 1)

$expectedDetailText = "Detailed instructions on ESP setup for steps 2-4 can be found in the DEPA-INSTRUCT.docx word document. "
if ($detailContent.Text -ne $expectedDetailText) {
    throw "Unexpected detail text: $($detailContent.Text)"
}

$detailXml = @'
<w:p>
<w:r><w:t xml:space="preserve">Detailed instructions on ESP setup for steps 2-4 can be found in the DEPA-INSTRUCT.docx word document. </w:t></w:r>
</w:p>
'@

$detailContent.InsertXML((Get-PkgXml $detailXml))

Write-Output "Done."
